# Apply updated cryptocurrency price/volume figures to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value, and whether the value must be
# forced to stay text (prefixed with a leading apostrophe) because it would
# otherwise be auto-recognized by Excel as a number and lose formatting such
# as trailing zeros (e.g. "1.0000" -> 1).
$updates = @(
    @{ Cell = 'D2'; Value = '30.365.24'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -0.59%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.871.11'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D4'; Value = '0.9998'; ForceText = $true }
    @{ Cell = 'E4'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '243.92'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.62%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '1.0000'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  -1.31%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.2875'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -1.15%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '0.06450'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -0.76%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '22.07'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +0.88%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.07759'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '1.868.77'; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '95.92'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  -0.23%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '0.7232'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -2.20%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '5.130'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -0.76%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '278.93'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +1.58%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '30.352.56'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -0.81%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '12.97'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -1.84%  '; ForceText = $false }
    @{ Cell = 'E19'; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.000007497'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '2.112.99'; ForceText = $false }
    @{ Cell = 'E21'; Value = '  -0.25%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '0.9997'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '5.225'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '6.225'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +0.83%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '163.13'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -1.09%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '9.050'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -1.51%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '18.67'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -0.68%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '1.872'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -2.00%  '; ForceText = $false }
    @{ Cell = 'E29'; Value = '  -1.22%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '0.09615'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -2.38%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '1.476'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -1.39%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '4.201'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -1.31%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '4.104'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +0.51%  '; ForceText = $false }
    @{ Cell = 'E34'; Value = '  +0.14%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.118'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '0.6884'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -0.90%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '2.715'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.01873'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.94%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '2.810'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +1.90%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '6.197'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -1.41%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '74.19'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +1.25%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.4224'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +0.42%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '1.930'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -2.49%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.8289'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -0.73%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '100.78'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -0.90%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '9.540'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +1.73%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '35.27'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '6.938'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -0.48%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '904.39'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -0.95%  '; ForceText = $false }
    @{ Cell = 'E51'; Value = '  +0.88%  '; ForceText = $false }
)

foreach ($update in $updates) {
    $text = $update.Value
    if ($update.ForceText) {
        $text = "'" + $text
    }
    $ws.Range($update.Cell).Value = $text
}
